# Apply updated Price (D) and Volume(1h) (E) values for the symbol list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format first so numeric-looking / percent strings are
# preserved verbatim as text (matching the original inline-string cells)
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "330.05"
$ws.Range("E2").Value = "0.37%"
$ws.Range("D3").Value = "45.42"
$ws.Range("E3").Value = "2.78%"
$ws.Range("D4").Value = "5.476"
$ws.Range("E4").Value = "0.06%"
$ws.Range("D5").Value = "0.08459"
$ws.Range("E5").Value = "4.89%"
$ws.Range("D6").Value = "2.049"
$ws.Range("E6").Value = "0.03%"
$ws.Range("D7").Value = "0.9805"
$ws.Range("E7").Value = "2.79%"
$ws.Range("D8").Value = "2.541"
$ws.Range("E8").Value = "-3.09%"
$ws.Range("D9").Value = "0.1144"
$ws.Range("E9").Value = "1.36%"
$ws.Range("D10").Value = "0.1909"
$ws.Range("E10").Value = "1.63%"
$ws.Range("D11").Value = "9.450"
$ws.Range("E11").Value = "-7.92%"
$ws.Range("D12").Value = "0.09681"
$ws.Range("E12").Value = "-2.87%"
$ws.Range("D13").Value = "0.04713"
$ws.Range("E13").Value = "-1.36%"
$ws.Range("D14").Value = "0.1058"
$ws.Range("E14").Value = "-0.39%"
$ws.Range("D15").Value = "0.001301"
$ws.Range("E15").Value = "2.44%"
$ws.Range("D16").Value = "0.04190"
$ws.Range("E16").Value = "2.39%"
$ws.Range("D17").Value = "0.005998"
$ws.Range("E17").Value = "2.64%"
$ws.Range("D18").Value = "3.387"
$ws.Range("E18").Value = "0.31%"
$ws.Range("D19").Value = "4.456"
$ws.Range("E19").Value = "0.85%"
$ws.Range("D20").Value = "0.3353"
$ws.Range("E20").Value = "-1.66%"
$ws.Range("D21").Value = "0.1359"
$ws.Range("E21").Value = "-3.01%"
$ws.Range("D22").Value = "0.2560"
$ws.Range("E22").Value = "-0.79%"
$ws.Range("D23").Value = "0.001303"
$ws.Range("E23").Value = "-0.24%"
$ws.Range("D24").Value = "0.004438"
$ws.Range("E24").Value = "2.14%"
$ws.Range("D25").Value = "0.0001303"
$ws.Range("E25").Value = "4.18%"
$ws.Range("D38").Value = "0.02709"
$ws.Range("E38").Value = "5.14%"
$ws.Range("D39").Value = "0.05679"
$ws.Range("E39").Value = "0.09%"
$ws.Range("D40").Value = "0.007821"
$ws.Range("E40").Value = "3.28%"
$ws.Range("D41").Value = "0.1426"
$ws.Range("E41").Value = "1.94%"
$ws.Range("D42").Value = "0.007478"
$ws.Range("E42").Value = "1.67%"
$ws.Range("D43").Value = "0.002121"
$ws.Range("E43").Value = "5.55%"
$ws.Range("D44").Value = "0.007904"
$ws.Range("E44").Value = "-7.08%"
$ws.Range("D45").Value = "0.3389"
$ws.Range("D46").Value = "0.00006953"
$ws.Range("E46").Value = "-0.89%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.20%"
$ws.Range("E48").Value = "0.16%"
$ws.Range("D49").Value = "0.003475"
$ws.Range("E49").Value = "-0.96%"
$ws.Range("D50").Value = "0.003550"
$ws.Range("E50").Value = "1.39%"
$ws.Range("D51").Value = "0.00002106"
$ws.Range("E51").Value = "0.20%"

# Restore default (General) formatting now that the text values are committed
$ws.Range("D2:E51").ClearFormats()
